$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data row: add "1372 / Bayer U.S. / CJ00617098" below the existing table,
# matching the formatting of the row above it (row 13).
$ws.Range("A13:C13").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)  # xlPasteFormats - keeps fill/border/font/number format

$ws.Range("A14").Value = "1372"
$ws.Range("B14").Value = "Bayer U.S."
$ws.Range("C14").Value = "CJ00617098"

# Match the row height used by the rest of the table.
$ws.Rows.Item(14).RowHeight = 21

# Reflect the saved selection/scroll state.
[void]$ws.Range("J5").Select()
